$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price (column D) and 1h volume change (column E) values,
# refreshed by the "Updated symbol list" GitHub Actions job.
$updates = @{
    "D2" = "291.41"
    "E2" = "-3.70%"
    "D3" = "30.89"
    "E3" = "-3.93%"
    "D4" = "4.874"
    "E4" = "-2.39%"
    "D5" = "0.07171"
    "E5" = "-9.35%"
    "E6" = "-2.54%"
    "D7" = "1.752"
    "E7" = "-16.66%"
    "D8" = "3.762"
    "E8" = "-1.20%"
    "D9" = "0.9002"
    "E9" = "-2.86%"
    "D10" = "0.1647"
    "E10" = "-6.17%"
    "D11" = "0.07420"
    "E11" = "-7.17%"
    "D12" = "0.08005"
    "E12" = "-8.57%"
    "D13" = "0.03048"
    "E13" = "-3.80%"
    "D14" = "0.09974"
    "E14" = "-0.72%"
    "D15" = "0.001503"
    "E15" = "-0.79%"
    "D16" = "0.005661"
    "E16" = "-5.82%"
    "D17" = "3.462"
    "E17" = "-0.12%"
    "D18" = "2.111"
    "E18" = "-7.38%"
    "D19" = "0.3292"
    "E19" = "0.16%"
    "D20" = "0.1305"
    "E20" = "1.12%"
    "E21" = "4.78%"
    "D22" = "0.2010"
    "E22" = "12.15%"
    "D23" = "0.04483"
    "E23" = "-2.78%"
    "D24" = "0.001219"
    "E24" = "-1.46%"
    "D25" = "0.004016"
    "E25" = "-10.54%"
    "D26" = "0.0001257"
    "E26" = "0.44%"
    "D39" = "0.01635"
    "E39" = "-5.91%"
    "E40" = "-9.72%"
    "D41" = "0.007411"
    "E41" = "1.11%"
    "D42" = "0.1314"
    "E42" = "-3.87%"
    "E43" = "-13.30%"
    "D44" = "0.01023"
    "E44" = "-7.71%"
    "D45" = "0.00005746"
    "E45" = "-4.50%"
    "E46" = "0.51%"
    "D47" = "2.178"
    "E47" = "164.48%"
    "D48" = "0.003016"
    "D49" = "0.00002111"
    "E49" = "0.51%"
    "E50" = "0.51%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = $origStyle
}
